$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Insert a new product row before the current row 9 ("ديتول صغير"),
# shifting it (and everything after it) down by one row. This makes
# room for the new "LOLAWEST 3GM 6 SACHETS" product line.
# ------------------------------------------------------------------
$ws.Rows("9:9").Insert()

# The plain Insert() does not carry over cell-level formatting
# (e.g. border), so clone it from the row above (row 8, the same
# kind of product row) one cell at a time - this reuses the existing
# style records instead of registering new duplicate ones.
$cols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q")
foreach ($col in $cols) {
    $ws.Range("${col}8").Copy()
    $ws.Range("${col}9").PasteSpecial(-4122)
}
$excel.CutCopyMode = $false

# Recreate the merges that Insert() doesn't add automatically for the
# freshly inserted row.
$ws.Range("A9:B9").Merge()
$ws.Range("C9:G9").Merge()
$ws.Range("H9:K9").Merge()
$ws.Range("L9:M9").Merge()
$ws.Range("N9:O9").Merge()

# ------------------------------------------------------------------
# Fill in the new product row (row 9). L9 and P9 use number-looking
# formats (165 / 0.00) so force them to stay text the way the source
# file stores them: switch to a plain text format, assign, then
# restore the original display format (paste-format only, value/type
# untouched).
# ------------------------------------------------------------------
$ws.Range("A9").Value = 3
$ws.Range("C9").Value = "LOLAWEST 3GM 6 SACHETS"
$ws.Range("H9").Value = "1:0"

$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "1"
$ws.Range("L8").Copy()
$ws.Range("L9").PasteSpecial(-4122)

$ws.Range("N9").Value = "96.00"

$ws.Range("P9").NumberFormat = "@"
$ws.Range("P9").Value = "96.0000"
$ws.Range("P8").Copy()
$ws.Range("P9").PasteSpecial(-4122)

$ws.Range("Q9").Value = "1:0"
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# Renumber the two product rows that were pushed down.
# ------------------------------------------------------------------
$ws.Range("A10").Value = 4
$ws.Range("A11").Value = 5

# ------------------------------------------------------------------
# Update the running total (old 171 -> new 267 after adding 96.00).
# ------------------------------------------------------------------
$ws.Range("P12").Value = 267

# ------------------------------------------------------------------
# Refresh the generated-at timestamp in the footer.
# ------------------------------------------------------------------
$ws.Range("A13").Value = "Monday, 11 August, 2025 9:56 AM"
